# Add a new "comet_score" worksheet right after "composite_score",
# matching the shape/style of the other per-metric sheets
# (bertscore_f1 / bleu / chrf / labse_similarity): row 1 holds the
# language headers in B1:I1, column A (rows 2-5) holds the model
# names, and B2:I5 holds the per-model / per-language COMET scores.

$wb = $excel.ActiveWorkbook

$after = $wb.Worksheets.Item("composite_score")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$ws.Name = "comet_score"

# match the page margins used by the sibling metric sheets (0.75in / 1in / 0.5in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# --- header row (languages) ---
$ws.Range("B1").Value = "Arabic"
$ws.Range("C1").Value = "Chinese (Simplified)"
$ws.Range("D1").Value = "Haitian Creole"
$ws.Range("E1").Value = "Korean"
$ws.Range("F1").Value = "Russian"
$ws.Range("G1").Value = "Spanish"
$ws.Range("H1").Value = "Tagalog/Filipino"
$ws.Range("I1").Value = "Vietnamese"

# --- row labels (models) ---
$ws.Range("A2").Value = "claude-opus-4.5"
$ws.Range("A3").Value = "gemini-3-pro"
$ws.Range("A4").Value = "gpt-5.1"
$ws.Range("A5").Value = "kimi-k2"

# --- data: comet_score per model x language ---
# claude-opus-4.5
$ws.Range("B2").Value = 0.9107455064853033
$ws.Range("C2").Value = 0.9086133142312368
$ws.Range("D2").Value = 0.9111654659112295
$ws.Range("E2").Value = 0.90642083187898
$ws.Range("F2").Value = 0.9081218242645264
$ws.Range("G2").Value = 0.9156241019566854
$ws.Range("H2").Value = 0.9132047444581985
$ws.Range("I2").Value = 0.9107471654812495

# gemini-3-pro
$ws.Range("B3").Value = 0.896774227420489
$ws.Range("C3").Value = 0.8985321124394735
$ws.Range("D3").Value = 0.8927175501982371
$ws.Range("E3").Value = 0.891526406009992
$ws.Range("F3").Value = 0.8960749258597692
$ws.Range("G3").Value = 0.9106310854355494
$ws.Range("H3").Value = 0.8972395757834116
$ws.Range("I3").Value = 0.8977503925561905

# gpt-5.1
$ws.Range("B4").Value = 0.8983265161514282
$ws.Range("C4").Value = 0.9020032833019892
$ws.Range("D4").Value = 0.9060116608937582
$ws.Range("E4").Value = 0.8954385370016098
$ws.Range("F4").Value = 0.9036827882130941
$ws.Range("G4").Value = 0.8994199832280477
$ws.Range("H4").Value = 0.9104017019271851
$ws.Range("I4").Value = 0.9089181621869405

# kimi-k2
$ws.Range("B5").Value = 0.8885662903388342
$ws.Range("C5").Value = 0.8956655164559683
$ws.Range("D5").Value = 0.853369931379954
$ws.Range("E5").Value = 0.881712332367897
$ws.Range("F5").Value = 0.8639207283655802
$ws.Range("G5").Value = 0.8981739282608032
$ws.Range("H5").Value = 0.8935400992631912
$ws.Range("I5").Value = 0.8768882602453232

# --- styling: bold + centered + thin border on header row and row labels,
#     matching the look of the other metric sheets ---
$headerRow = $ws.Range("B1:I1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1

$labelCol = $ws.Range("A2:A5")
$labelCol.Font.Bold = $true
$labelCol.HorizontalAlignment = -4108
$labelCol.VerticalAlignment = -4160
$labelCol.Borders.LineStyle = 1

[void]$ws.Range("A1").Select()
